$wb = $excel.ActiveWorkbook

# --- Sheet2 ("Sheet2" tab): update the small results table A4:C9 ---
$ws3 = $wb.Worksheets.Item("Sheet2")

# Shift the remaining two rows of data up, clearing the rest.
$ws3.Range("A4").Value = "DIP B1"
$ws3.Range("B4").Value = 3.6890645586297759
$ws3.Range("C4").Value = 0.12440240703059115

$ws3.Range("A5").Value = "DIP G2"
$ws3.Range("B5").Value = 8.4175084175084187
$ws3.Range("C5").Value = 0.734146272424678

$ws3.Range("A6:C6").ClearContents()
$ws3.Range("A7:C7").ClearContents()
$ws3.Range("A8:C8").ClearContents()
$ws3.Range("A9:C9").ClearContents()

# Update selection on the Sheet2 tab
$ws3.Range("A1:C1").Select()

# --- Sheet1 tab: move selection ---
$ws2 = $wb.Worksheets.Item("Sheet1")
$ws2.Range("E3").Select()

# --- Plate 1 - Sheet1 tab: scroll the view ---
$ws1 = $wb.Worksheets.Item("Plate 1 - Sheet1")
$excel.ActiveWindow.ScrollRow = 32

# Re-select Sheet2 tab (it was tabSelected in the original file)
$ws3.Select()
